$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: single value tweak ---
$ws.Range("B2").Value = 2

# --- Row 3: turn several literals into formulas ---
$ws.Range("A3").Formula = "=A8"
$ws.Range("C3").Formula = "=C5*2"
$ws.Range("D3:E3").Formula = "=D5*2"
$ws.Range("F3").Formula = "=7"

# --- Row 8: previously-empty row, now filled in (copy of row 5's pattern,
#     but referencing row 6 two rows above it, same as row 5 references row 3) ---
$ws.Range("A8").Value = 28
$ws.Range("B8").Value = 28
$ws.Range("C8").Value = 8
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 7
$ws.Range("G8:T8").Formula = "=G6*2"

# Match the formatting used on the rest of row 5 (style index carries no
# visible difference from default, but keep cells consistently formatted).
$ws.Range("A8:T8").Style = $ws.Range("A5:T5").Style

# --- Selection, as left by the author ---
$ws.Range("E6").Select()
